# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates column G ("K") values on the active sheet (rows 2-19)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 6
    3  = 3
    4  = 1
    5  = 3
    6  = 4
    7  = 2
    8  = 4
    9  = 3
    10 = 2
    11 = 8
    12 = 4
    13 = 6
    14 = 2
    15 = 2
    16 = 1
    17 = 1
    18 = 3
    19 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
